$d = $word.ActiveDocument

# --- Rewrite existing paragraphs to add proofErr (spell-check) markup ---
$r2 = $d.Paragraphs(2).Range
$r2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t xml:space="preserve">VM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Scale</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> set</w:t></w:r></w:p>')

$r4 = $d.Paragraphs(4).Range
$r4.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t xml:space="preserve">Execução e orquestração de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>containeres</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$r5 = $d.Paragraphs(5).Range
$r5.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Azure container </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>instances</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: útil para executar uma quantidade pequena e pré-definida de instâncias de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>contâiner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r></w:p>')

$r6 = $d.Paragraphs(6).Range
$r6.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Azue</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> servisse </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fabric</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>custom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Microsoft;</w:t></w:r></w:p>')

$r7 = $d.Paragraphs(7).Range
$r7.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">ARO: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Openshift</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> PaaS para Azure;</w:t></w:r></w:p>')

$r8 = $d.Paragraphs(8).Range
$r8.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">AKS: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Kubernetes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> PaaS para Azure;</w:t></w:r></w:p>')

$r9 = $d.Paragraphs(9).Range
$r9.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Azure container apps: similar ao “Azure container </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>instances</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">”, porém possui uma certa escalabilidade, sem que haja uma orquestração do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Kubernetes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r></w:p>')

$r10 = $d.Paragraphs(10).Range
$r10.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Chaos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>studio</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$r11 = $d.Paragraphs(11).Range
$r11.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Elemento similar ao “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Chaos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>monkey</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” para criar testes de caos na solução de aplicações Azure</w:t></w:r></w:p>')

$r12 = $d.Paragraphs(12).Range
$r12.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t xml:space="preserve">Data </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>storage</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$r13 = $d.Paragraphs(13).Range
$r13.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: solução simples para armazenamento de objetos chave/valor</w:t></w:r></w:p>')

$r14 = $d.Paragraphs(14).Range
$r14.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Queue</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$r15 = $d.Paragraphs(15).Range
$r15.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">File </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>share</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$r16 = $d.Paragraphs(16).Range
$r16.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Todos os elementos do “Data </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>storage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” possuem um SDK compatível com a tecnologia escolhida (Java, .net, etc.)</w:t></w:r></w:p>')

$r19 = $d.Paragraphs(19).Range
$r19.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Azure </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Powershell</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$r21 = $d.Paragraphs(21).Range
$r21.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Azure BICEP / ARM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Template</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$r22 = $d.Paragraphs(22).Range
$r22.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Terraform</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IaC</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>')

$r25 = $d.Paragraphs(25).Range
$r25.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Policies podem ser definidas em dois níveis: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tenant</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> e </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>subscription</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# --- Append new paragraphs at the end of the document ---
$endRange0 = $d.Content
$endRange0.Collapse(0)
$endRange0.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t xml:space="preserve">Azure </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bastion</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$endRange1 = $d.Content
$endRange1.Collapse(0)
$endRange1.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Jumpserver</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jumpbox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) para proteção de recursos acessados por hosts de origens inseguras</w:t></w:r></w:p>')
$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br w:type="page"/></w:r></w:p>')
$endRange3 = $d.Content
$endRange3.Collapse(0)
$endRange3.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Objetos de infra criados</w:t></w:r></w:p>')
$endRange4 = $d.Content
$endRange4.Collapse(0)
$endRange4.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Tipo: virtual </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>machine</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$endRange5 = $d.Content
$endRange5.Collapse(0)
$endRange5.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Usuário: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rctbatista</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$endRange6 = $d.Content
$endRange6.Collapse(0)
$endRange6.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Senha: </w:t></w:r><w:r><w:t>rctbatista</w:t></w:r><w:r><w:t>@D@ttebayo5</w:t></w:r></w:p>')
$endRange7 = $d.Content
$endRange7.Collapse(0)
$endRange7.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br w:type="page"/></w:r></w:p>')
$endRange8 = $d.Content
$endRange8.Collapse(0)
$endRange8.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Link de avalições de aula</w:t></w:r></w:p>')
$endRange9 = $d.Content
$endRange9.Collapse(0)
$endRange9.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Aula 1: </w:t></w:r><w:r><w:t>https://mastertech-tech.typeform.com/to/X9EhFux1</w:t></w:r></w:p>')

Write-Output "done"
